$d = $word.ActiveDocument

# Replace all occurrences of "Constel·lació de bessons" with
# "Constel·lació de Bessons" (capitalize "Bessons") throughout the document.
$d.Content.Find.Execute("Constel·lació de bessons", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Constel·lació de Bessons", 2)
